# Apply cryptocurrency price/volume updates scraped on Sat Oct 28 18:41:52 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.179.49"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.784.91"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.73"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").Value = "  +2.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0688"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.041.87"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.44%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.01"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.785.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.109.66"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.622"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.24"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.19%  "

$ws.Range("E20").Value = "  +1.29%  "

$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.90%  "

$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("E24").Value = "  +0.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.59%  "

$ws.Range("E26").Value = "  +3.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.87%  "

$ws.Range("E28").Value = "  +2.22%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0519"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.08%  "

$ws.Range("E33").Value = "  +4.54%  "

$ws.Range("E34").Value = "  +1.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.446.30"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.32%  "

$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.43"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.66%  "

$ws.Range("E38").Value = "  +4.60%  "

$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.21"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.87%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("E42").Value = "  +2.95%  "

$ws.Range("E43").Value = "  +1.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.50"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.07"
$ws.Range("D45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0508"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.13%  "

$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.943.93"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.74"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("E51").Value = "  -0.01%  "
